# Reproduce the "add captions in run file" commit:
#  - add a new worksheet "add" (after "withdraw", i.e. as sheet 5)
#  - put the usual case-table header row in it plus two short data rows,
#    the second one captioned "添加成功" (new shared string)
#  - leave the other sheets' last-used-cell selections the way the
#    author left them while doing this edit, and finish with the new
#    "add" sheet active/selected (matches the new activeTab="4").

$wb = $excel.ActiveWorkbook

# --- register sheet: author was last at B4 ------------------------------
$wsRegister = $wb.Worksheets.Item("register")
$wsRegister.Activate()
$wsRegister.Range("B4").Select()

# --- recharge sheet: author was last at E2 -------------------------------
$wsRecharge = $wb.Worksheets.Item("recharge")
$wsRecharge.Activate()
$wsRecharge.Range("E2").Select()

# --- withdraw sheet: author selected the whole header row ----------------
$wsWithdraw = $wb.Worksheets.Item("withdraw")
$wsWithdraw.Activate()
$wsWithdraw.Rows("1:1").Select()

# --- new "add" sheet, inserted right after "withdraw" --------------------
$newSheet = $wb.Worksheets.Add($null, $wsWithdraw)
$newSheet.Name = "add"

# Reuse the same header look (fill/border) as the other case sheets by
# copying the formats of the "withdraw" header row onto the new one.
$wsWithdraw.Range("A1:G1").Copy()
$newSheet.Range("A1:G1").PasteSpecial(-4122)

$newSheet.Range("A1").Value = "case_id"
$newSheet.Range("B1").Value = "title"
$newSheet.Range("C1").Value = "url"
$newSheet.Range("D1").Value = "method"
$newSheet.Range("E1").Value = "data"
$newSheet.Range("F1").Value = "expected"
$newSheet.Range("G1").Value = "sql"

$newSheet.Range("A2").Value = 1
$newSheet.Range("B2").Value = "添加成功"
$newSheet.Range("A3").Value = 2

# Make the new sheet the active tab, with the caret left at N9.
$newSheet.Activate()
$newSheet.Range("N9").Select()
